$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rnaDate (col A) and s1cDNADate (col D) values for this run.
# Rows 2-10 move from 12.06.11 -> 12.07.11 (A) / 12.06.11 -> 12.08.11 (D)
# Rows 11-18 move from 12.07.11 -> 12.08.11 (A) / 12.07.11 -> 12.09.11 (D)
# Force the cells to stay text (matching the original "General" shared-string
# cells) instead of letting the date-looking strings auto-convert to dates.

$colA = $ws.Range("A2:A18")
$colD = $ws.Range("D2:D18")
$colA.NumberFormat = "@"
$colD.NumberFormat = "@"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = "12.07.11"
    $ws.Cells.Item($r, 4).Value = "12.08.11"
}

for ($r = 11; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = "12.08.11"
    $ws.Cells.Item($r, 4).Value = "12.09.11"
}

$colA.Style = "Normal"
$colD.Style = "Normal"

# Update the selection to match the new active cell/range.
$ws.Range("D12:D18").Select()
